$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 5053
$ws.Cells.Item(2, 2).Value = 45919
$ws.Cells.Item(3, 1).Value = 5063
$ws.Cells.Item(3, 2).Value = 45919.01041666666
$ws.Cells.Item(4, 1).Value = 5051
$ws.Cells.Item(4, 2).Value = 45919.02083333334
$ws.Cells.Item(5, 1).Value = 5031
$ws.Cells.Item(5, 2).Value = 45919.03125
$ws.Cells.Item(6, 1).Value = 4981
$ws.Cells.Item(6, 2).Value = 45919.04166666666
$ws.Cells.Item(7, 1).Value = 5038
$ws.Cells.Item(7, 2).Value = 45919.05208333334
$ws.Cells.Item(8, 1).Value = 4972
$ws.Cells.Item(8, 2).Value = 45919.0625
$ws.Cells.Item(9, 1).Value = 4985
$ws.Cells.Item(9, 2).Value = 45919.07291666666
$ws.Cells.Item(10, 1).Value = 4933
$ws.Cells.Item(10, 2).Value = 45919.08333333334
$ws.Cells.Item(11, 1).Value = 4949
$ws.Cells.Item(11, 2).Value = 45919.09375
$ws.Cells.Item(12, 1).Value = 4906
$ws.Cells.Item(12, 2).Value = 45919.10416666666
$ws.Cells.Item(13, 1).Value = 4922
$ws.Cells.Item(13, 2).Value = 45919.11458333334
$ws.Cells.Item(14, 1).Value = 4942
$ws.Cells.Item(14, 2).Value = 45919.125
$ws.Cells.Item(15, 1).Value = 5011
$ws.Cells.Item(15, 2).Value = 45919.13541666666
$ws.Cells.Item(16, 1).Value = 4974
$ws.Cells.Item(16, 2).Value = 45919.14583333334
$ws.Cells.Item(17, 1).Value = 4991
$ws.Cells.Item(17, 2).Value = 45919.15625
$ws.Cells.Item(18, 1).Value = 5117
$ws.Cells.Item(18, 2).Value = 45919.16666666666
$ws.Cells.Item(19, 1).Value = 5198
$ws.Cells.Item(19, 2).Value = 45919.17708333334
$ws.Cells.Item(20, 1).Value = 5168
$ws.Cells.Item(20, 2).Value = 45919.1875
$ws.Cells.Item(21, 1).Value = 5262
$ws.Cells.Item(21, 2).Value = 45919.19791666666
$ws.Cells.Item(22, 1).Value = 5484
$ws.Cells.Item(22, 2).Value = 45919.20833333334
$ws.Cells.Item(23, 1).Value = 5587
$ws.Cells.Item(23, 2).Value = 45919.21875
$ws.Cells.Item(24, 1).Value = 5640
$ws.Cells.Item(24, 2).Value = 45919.22916666666
$ws.Cells.Item(25, 1).Value = 5709
$ws.Cells.Item(25, 2).Value = 45919.23958333334
$ws.Cells.Item(26, 1).Value = 5941
$ws.Cells.Item(26, 2).Value = 45919.25
$ws.Cells.Item(27, 1).Value = 6121
$ws.Cells.Item(27, 2).Value = 45919.26041666666
$ws.Cells.Item(28, 1).Value = 6129
$ws.Cells.Item(28, 2).Value = 45919.27083333334
$ws.Cells.Item(29, 1).Value = 6160
$ws.Cells.Item(29, 2).Value = 45919.28125
$ws.Cells.Item(30, 1).Value = 6230
$ws.Cells.Item(30, 2).Value = 45919.29166666666
$ws.Cells.Item(31, 1).Value = 6284
$ws.Cells.Item(31, 2).Value = 45919.30208333334
$ws.Cells.Item(32, 1).Value = 6235
$ws.Cells.Item(32, 2).Value = 45919.3125
$ws.Cells.Item(33, 1).Value = 6082
$ws.Cells.Item(33, 2).Value = 45919.32291666666
$ws.Cells.Item(34, 1).Value = 5983
$ws.Cells.Item(34, 2).Value = 45919.33333333334
$ws.Cells.Item(35, 1).Value = 5905
$ws.Cells.Item(35, 2).Value = 45919.34375
$ws.Cells.Item(36, 1).Value = 5815
$ws.Cells.Item(36, 2).Value = 45919.35416666666
$ws.Cells.Item(37, 1).Value = 5680
$ws.Cells.Item(37, 2).Value = 45919.36458333334
$ws.Cells.Item(38, 1).Value = 5577
$ws.Cells.Item(38, 2).Value = 45919.375
$ws.Cells.Item(39, 1).Value = 5473
$ws.Cells.Item(39, 2).Value = 45919.38541666666
$ws.Cells.Item(40, 1).Value = 5343
$ws.Cells.Item(40, 2).Value = 45919.39583333334
$ws.Cells.Item(41, 1).Value = 5243
$ws.Cells.Item(41, 2).Value = 45919.40625
$ws.Cells.Item(42, 1).Value = 5105
$ws.Cells.Item(42, 2).Value = 45919.41666666666
